$wb = $excel.ActiveWorkbook

# "Generate Report for Archive" — the localization-status report was
# regenerated; 1ed60ced-...md and 64459f6e-...md moved from
# "Ready for handoff" to "In Translation" (a8348097-...md stays
# "Ready for handoff"). Update every sheet that surfaces the status text.

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("E4").Value = "In Translation"
$overview.Range("F4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"
